$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 1359
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = 842
$ws.Range("F11").Value = 842
$ws.Range("F12").Value = 129
$ws.Range("F13").Value = 2929
$ws.Range("F14").Value = 395
$ws.Range("F15").Value = 886
$ws.Range("F16").Value = 1138
$ws.Range("F17").Value = 595
$ws.Range("F18").Value = 338
$ws.Range("F19").Value = 71
$ws.Range("F20").Value = 1673
$ws.Range("F21").Value = 348
$ws.Range("F22").Value = 1269
$ws.Range("F23").Value = 218
$ws.Range("F24").Value = 602
$ws.Range("F26").Value = 1077
$ws.Range("F27").Value = 1534
$ws.Range("F28").Value = 1476
$ws.Range("F29").Value = 1345
$ws.Range("F30").Value = 356
$ws.Range("F33").Value = 158
$ws.Range("F34").Value = 977
$ws.Range("F36").Value = 1859
$ws.Range("F37").Value = 488
$ws.Range("F38").Value = 1056
$ws.Range("F40").Value = 22
$ws.Range("F41").Value = 2303
$ws.Range("F44").Value = 2808
$ws.Range("F47").Value = 648
$ws.Range("F49").Value = 12
$ws.Range("F50").Value = 34

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 114671
$ws.Range("F17").Value = 74
$ws.Range("F18").Value = 74
$ws.Range("F22").Value = 288
$ws.Range("F30").Value = 48
$ws.Range("F38").Value = 163
$ws.Range("F40").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 3033
$ws.Range("F6").Value = 4858
$ws.Range("F12").Value = 648
$ws.Range("F13").Value = 1374
$ws.Range("F14").Value = 390
$ws.Range("F15").Value = 1292

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 4858
$ws.Range("F9").Value = 648
$ws.Range("F10").Value = 1374
$ws.Range("F13").Value = 1359
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 843
$ws.Range("F16").Value = 843
$ws.Range("F17").Value = 2929
$ws.Range("F19").Value = 395
$ws.Range("F20").Value = 886
$ws.Range("F21").Value = 1139
$ws.Range("F22").Value = 595
$ws.Range("F23").Value = 338
$ws.Range("F24").Value = 1674
$ws.Range("F26").Value = 348
$ws.Range("F28").Value = 602
$ws.Range("F30").Value = 1077
$ws.Range("F31").Value = 1534
$ws.Range("F32").Value = 1476
$ws.Range("F33").Value = 1345
$ws.Range("F34").Value = 356
$ws.Range("F35").Value = 74
$ws.Range("F38").Value = 977
$ws.Range("F40").Value = 1859
$ws.Range("F42").Value = 1056
$ws.Range("F45").Value = 2303
$ws.Range("F48").Value = 0
$ws.Range("F50").Value = 648
$ws.Range("F52").Value = 2

Write-Output "Applied all F-column updates"